$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Items.ItemType enum definition + first item "Sword"
$ws.Range("B7").Value = "Items.ItemType"
$ws.Range("D7").Value = $true
$ws.Range("H7").Value = "Sword"

# Row 8: Items.Rarity enum definition + first item "Common"
$ws.Range("B8").Value = "Items.Rarity"
$ws.Range("D8").Value = $true
$ws.Range("H8").Value = "Common"

# Rows 9-12: remaining Items.Rarity values
$ws.Range("H9").Value = "Rare"
$ws.Range("H10").Value = "Magic"
$ws.Range("H11").Value = "Epic"
$ws.Range("H12").Value = "Legendary"

# Match the author's final selection state
$ws.Range("H12").Select()
